# Update BTCUSDTWCHG sheet: correct the last existing row (262) and append
# the new weekly OHLC rows (263-271) that were added in the source export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up row 262 (High / Low / Close changed) -----------------------
$ws.Cells.Item(262, 3).Value = 25047.56
$ws.Cells.Item(262, 4).Value = 22664.69
$ws.Cells.Item(262, 5).Value = 24305.24

# --- Append the new rows --------------------------------------------------
# Columns: Open time (A, date serial, same number format as the rest of the
# column), Open (B), High (C), Low (D), Close (E)
$newRows = @(
    @(263, 44788, 24305.25, 25211.32, 20761.90, 21515.61),
    @(264, 44795, 21516.70, 21900.00, 19520.00, 19555.61),
    @(265, 44802, 19555.61, 20576.25, 19540.00, 20000.30),
    @(266, 44809, 20000.30, 21860.00, 18510.77, 21826.87),
    @(267, 44816, 21826.87, 22799.00, 19320.01, 19416.18),
    @(268, 44823, 19417.45, 19956.00, 18125.98, 18807.38),
    @(269, 44830, 18809.13, 20385.86, 18471.28, 19056.80),
    @(270, 44837, 19057.74, 20475.00, 18959.68, 19439.02),
    @(271, 44844, 19439.96, 19951.87, 18190.00, 19175.86)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $openTimeCell = $ws.Cells.Item($r, 1)
    $openTimeCell.Value = $row[1]
    # Match the date/time number format already used for column A
    $openTimeCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
